$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (and Aave/Frax row swap at 47-48).
# Values are written with a leading apostrophe to force text storage
# (matching the original inlineStr cells), then the style is reset to
# "Normal" so no stray quote-prefix formatting is left behind.

$ws.Range("D2").Value = "'25.911.77"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").Value = "'1.636.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.32%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'216.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.75%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'0.5069"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.01%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D8").Value = "'0.2584"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.90%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.06371"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.02%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'19.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.82%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07757"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.08%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'4.272"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.39%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'1.634.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.85%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.5511"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.96%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.0₅7737"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.16%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'64.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.26%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'25.898.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.44%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("E18").Value = "'  -0.30%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'4.452"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.21%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'195.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.04%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'9.919"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.50%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  +0.59%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'1.003"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.43%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'1.904"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.39%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'142.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.79%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.1241"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +5.70%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'6.841"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.59%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'15.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.38%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'1.247"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.77%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'0.04874"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.19%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'3.249"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.36%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'3.202"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.47%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'1.547"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.27%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'2.372"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.41%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'0.9070"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.40%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'2.572"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.28%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.5513"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.28%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'1.123.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.74%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.01559"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.32%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  -0.33%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'5.588"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.29%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.8049"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.64%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'97.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.92%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value = "'  -4.54%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'1.772.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.35%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.4459"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.64%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Value = "'Aave"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'54.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.01%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("B48").Value = "'Frax"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'0.9966"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.71%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.05146"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.48%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'7.527"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.91%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'1.004"
$ws.Range("D51").Style = "Normal"
